# This script updates the "want-to-go" visitor counts (column F) across the
# three sheets that mirror event rows: "展览" (Exhibitions), "演出" (Performances),
# and "全部类型" (All Types, a combined listing). These numbers are refreshed
# counters scraped from the source site, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsPerformance = $wb.Worksheets.Item("演出")
$wsAllTypes = $wb.Worksheets.Item("全部类型")

# -- 展览 (Exhibitions) --
$wsExhibition.Range("F2").Value = 277
$wsExhibition.Range("F4").Value = 285
$wsExhibition.Range("F6").Value = 70
$wsExhibition.Range("F8").Value = 2269
$wsExhibition.Range("F9").Value = 1518
$wsExhibition.Range("F13").Value = 2599
$wsExhibition.Range("F15").Value = 1443
$wsExhibition.Range("F16").Value = 6333
$wsExhibition.Range("F17").Value = 11
$wsExhibition.Range("F18").Value = 6202
$wsExhibition.Range("F20").Value = 2129
$wsExhibition.Range("F21").Value = 2987
$wsExhibition.Range("F22").Value = 3408
$wsExhibition.Range("F23").Value = 198
$wsExhibition.Range("F24").Value = 9
$wsExhibition.Range("F25").Value = 1680
$wsExhibition.Range("F26").Value = 58
$wsExhibition.Range("F30").Value = 16
$wsExhibition.Range("F31").Value = 347
$wsExhibition.Range("F32").Value = 1064
$wsExhibition.Range("F33").Value = 2283
$wsExhibition.Range("F35").Value = 139
$wsExhibition.Range("F37").Value = 861
$wsExhibition.Range("F38").Value = 176
$wsExhibition.Range("F39").Value = 411
$wsExhibition.Range("F40").Value = 477

# -- 演出 (Performances) --
$wsPerformance.Range("F3").Value = 107
$wsPerformance.Range("F15").Value = 3
$wsPerformance.Range("F17").Value = 19
$wsPerformance.Range("F19").Value = 90
$wsPerformance.Range("F22").Value = 56

# -- 全部类型 (All Types) --
$wsAllTypes.Range("F2").Value = 107
$wsAllTypes.Range("F5").Value = 277
$wsAllTypes.Range("F7").Value = 285
$wsAllTypes.Range("F10").Value = 70
$wsAllTypes.Range("F11").Value = 2269
$wsAllTypes.Range("F12").Value = 1518
$wsAllTypes.Range("F17").Value = 2599
$wsAllTypes.Range("F18").Value = 1443
$wsAllTypes.Range("F21").Value = 3
$wsAllTypes.Range("F23").Value = 6333
$wsAllTypes.Range("F24").Value = 11
$wsAllTypes.Range("F25").Value = 6202
$wsAllTypes.Range("F26").Value = 2129
$wsAllTypes.Range("F27").Value = 2987
$wsAllTypes.Range("F28").Value = 3408
$wsAllTypes.Range("F29").Value = 19
$wsAllTypes.Range("F30").Value = 198
$wsAllTypes.Range("F32").Value = 90
$wsAllTypes.Range("F33").Value = 1680
$wsAllTypes.Range("F39").Value = 16
$wsAllTypes.Range("F40").Value = 347
$wsAllTypes.Range("F41").Value = 56
$wsAllTypes.Range("F42").Value = 2284
$wsAllTypes.Range("F44").Value = 139
$wsAllTypes.Range("F46").Value = 861
$wsAllTypes.Range("F47").Value = 176
$wsAllTypes.Range("F48").Value = 411
$wsAllTypes.Range("F49").Value = 477

